$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 36426
$ws.Range("B2").Value = "Théo Dias"
$ws.Range("C2").Value = "P&D"
$ws.Range("D2").Value = "Outros"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45102
$ws.Range("G2").Value = 6641.16

# Row 3
$ws.Range("A3").Value = 25803
$ws.Range("B3").Value = "Giovanna das Neves"
$ws.Range("C3").Value = "Operacoes"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45101
$ws.Range("G3").Value = 6059.1

# Row 4
$ws.Range("A4").Value = 40077
$ws.Range("B4").Value = "Miguel Barbosa"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Viagem de negocios"
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45092
$ws.Range("G4").Value = 9642

# Row 5
$ws.Range("A5").Value = 2086
$ws.Range("B5").Value = "Igor Lima"
$ws.Range("C5").Value = "Juridico"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45084
$ws.Range("G5").Value = 2358.94

# Row 6
$ws.Range("A6").Value = 82207
$ws.Range("B6").Value = "Breno Nascimento"
$ws.Range("C6").Value = "TI"
$ws.Range("D6").Value = "Doenca"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45084
$ws.Range("G6").Value = 4430.35

# Row 7
$ws.Range("A7").Value = 44426
$ws.Range("B7").Value = "Bella da Conceição"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 45080
$ws.Range("G7").Value = 8231.24

# Row 8
$ws.Range("A8").Value = 39202
$ws.Range("B8").Value = "Sr. Miguel Moreira"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Doenca"
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 45079
$ws.Range("G8").Value = 2784.38

# Row 9
$ws.Range("A9").Value = 13316
$ws.Range("B9").Value = "Srta. Ana Lívia Souza"
$ws.Range("C9").Value = "Operacoes"
$ws.Range("D9").Value = "Doenca"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 45084
$ws.Range("G9").Value = 7248.5

# Row 10
$ws.Range("A10").Value = 43368
$ws.Range("B10").Value = "Marina Santos"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45080
$ws.Range("G10").Value = 6216.16

# Row 11
$ws.Range("A11").Value = 28576
$ws.Range("B11").Value = "Sr. Ryan Araújo"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45082
$ws.Range("G11").Value = 8003.91
